# Automatic map update (2025-09-27 09:16:16)
# Remove the case "-542" (Cramer 2141) entry, which shifts all subsequent
# rows up by one and shrinks the used range from A1:R40 to A1:R39.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

# Row 35 corresponds to case -542 / Cramer 2141, which has been removed.
# Deleting the entire row shifts rows 36:40 up to 35:39.
$ws.Rows.Item(35).Delete()
